$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: D11 "EMPTY" placeholder becomes a single space
$ws.Range("D11").Value = " "

# Update column L (反映内容) cells: replace old 外网编号 "电[YYYY]A-XXXXXXXX" format
# with the new plain numeric outer-network-number format (matching column A),
# and fill in the previously blank "反映内容：-" with its real content value.
$ws.Range("L2").Value = "外网编号：201939415931`n来点人姓名：张先生`n联系电话：1381188496`n主要内容：不开发票`n反映内容：5"
$ws.Range("L3").Value = "外网编号：201939515539`n来点人姓名：张先生`n联系电话：13910184444`n主要内容：税务政策不合理`n反映内容：5"
$ws.Range("L4").Value = "外网编号：201939488033`n来点人姓名：李先生`n联系电话：13910093936`n主要内容：华联超市物业不开发票“12350热线来电”`n反映内容：8"
$ws.Range("L5").Value = "外网编号：201939543342`n来点人姓名：胡女士`n联系电话：13562142093`n主要内容：映税务局迟迟没有回复`n反映内容：3"
$ws.Range("L6").Value = "外网编号：201954135923`n来点人姓名：朴女士`n联系电话：13010615197`n主要内容：不给销发票`n反映内容：5"
$ws.Range("L7").Value = "外网编号：201939562042`n来点人姓名：卢女士`n联系电话：13581956378`n主要内容：停车场收费不给开发票问题”12319热线来电“`n反映内容：2"
$ws.Range("L8").Value = "外网编号：201954233014`n来点人姓名：司先生`n联系电话：18513349999`n主要内容：不开发票`n反映内容：3"
$ws.Range("L9").Value = "外网编号：201954193509`n来点人姓名：朴先生`n联系电话：017110000019`n主要内容：望京小腰不开发票`n反映内容：3"
$ws.Range("L10").Value = "外网编号：201939612936`n来点人姓名：王先生`n联系电话：13811945536`n主要内容：老巷子家常菜商家未开发票`n反映内容：2"
$ws.Range("L11").Value = "外网编号：201954115023`n来点人姓名： `n联系电话：13011840840`n主要内容：偷税漏税问题`n反映内容：5"
$ws.Range("L12").Value = "外网编号：201939491932`n来点人姓名：刘女士`n联系电话：13880020880`n主要内容：超市不给开发票`n反映内容：1"
$ws.Range("L13").Value = "外网编号：201954191515`n来点人姓名：周先生`n联系电话：18852565888`n主要内容：停车场定额发票不盖章`n反映内容：1"
